# Week 15 log + Week 16 simulation update
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# YDS sheet: append newest-week numbers to the running per-play lists
# ---------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value2 = $wsYDS.Range("B2").Value2 + " 0 6 2 5 1 22 0 3 23 6 1 2 1 6 3 0 4 12"
$wsYDS.Range("C2").Value2 = $wsYDS.Range("C2").Value2 + " 10 0 7 6 -3 5 5 0 27 6 8 0 0 -3 1 0 1 4 5 3 0 6 0 2 8 1 16 5 0 2 1 2 -2 3"
$wsYDS.Range("B3").Value2 = $wsYDS.Range("B3").Value2 + " 6 -1 14 13 2 5 24 10 15 0 3 9 20 19 27 16 5 7 6 7 13 11 26"
$wsYDS.Range("C3").Value2 = $wsYDS.Range("C3").Value2 + " 4 5 15 5 13 5 13 5 37 5 15 14 8 22 16 6 5 8 11 -2"

# ---------------------------------------------------------------
# OFF sheet: Week totals updated
# ---------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value2 = 189
$wsOFF.Range("F2").Value2 = 40
$wsOFF.Range("G2").Value2 = 56
$wsOFF.Range("J2").Value2 = 29
$wsOFF.Range("N2").Value2 = 14
$wsOFF.Range("O2").Value2 = 21

$wsOFF.Range("B3").Value2 = 10
$wsOFF.Range("C3").Value2 = 143
$wsOFF.Range("E3").Value2 = 29
$wsOFF.Range("F3").Value2 = 106
$wsOFF.Range("G3").Value2 = 43
$wsOFF.Range("H3").Value2 = 26
$wsOFF.Range("I3").Value2 = 38
$wsOFF.Range("J3").Value2 = 46
$wsOFF.Range("L3").Value2 = 256
$wsOFF.Range("M3").Value2 = 185
$wsOFF.Range("Q3").Value2 = 533

# ---------------------------------------------------------------
# DEF sheet: Week totals updated
# ---------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("B2").Value2 = 7
$wsDEF.Range("C2").Value2 = 168
$wsDEF.Range("E2").Value2 = 9
$wsDEF.Range("F2").Value2 = 53
$wsDEF.Range("G2").Value2 = 39
$wsDEF.Range("I2").Value2 = 6
$wsDEF.Range("J2").Value2 = 22
$wsDEF.Range("N2").Value2 = 28
$wsDEF.Range("O2").Value2 = 26
$wsDEF.Range("P2").Value2 = 11

$wsDEF.Range("C3").Value2 = 138
$wsDEF.Range("E3").Value2 = 32
$wsDEF.Range("F3").Value2 = 80
$wsDEF.Range("G3").Value2 = 32
$wsDEF.Range("H3").Value2 = 25
$wsDEF.Range("I3").Value2 = 43
$wsDEF.Range("J3").Value2 = 51
$wsDEF.Range("L3").Value2 = 275
$wsDEF.Range("M3").Value2 = 179
$wsDEF.Range("Q3").Value2 = 494

# ---------------------------------------------------------------
# ST sheet: Week totals + per-kick / per-return lists
# ---------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value2 = 77
$wsST.Range("D2").Value2 = 45
$wsST.Range("H2").Value2 = 2
$wsST.Range("J2").Value2 = 211
$wsST.Range("K2").Value2 = 198

$wsST.Range("B3").Value2 = 61

$wsST.Range("B4").Value2 = $wsST.Range("B4").Value2 + " 64"
$wsST.Range("B5").Value2 = $wsST.Range("B5").Value2 + " 17"
$wsST.Range("B6").Value2 = $wsST.Range("B6").Value2 + " 17 31"
$wsST.Range("D3").Value2 = $wsST.Range("D3").Value2 + " 34 49 49"
$wsST.Range("D4").Value2 = $wsST.Range("D4").Value2 + " 0 10 14"
$wsST.Range("D5").Value2 = $wsST.Range("D5").Value2 + " 6 12 0"

# ---------------------------------------------------------------
# TURNS sheet: Week totals
# ---------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("B3").Value2 = 3
$wsTURNS.Range("D3").Value2 = 9
$wsTURNS.Range("E3").Value2 = 7

# ---------------------------------------------------------------
# PEN sheet: Week totals
# ---------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("B2").Value2 = 21
$wsPEN.Range("D4").Value2 = 6
